# Auto-generated Excel COM-interop script
# Adds skin disease data (Ringworm, Warts Molluscum, Serorrheic Keratoses, Nail Fungus)
# plus home-remedy text for Skin Cancer, Melanoma, Melanocytic Nevi and Atopic Dermatitis rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 in the original sheet is a fully blank spacer row; deleting it shifts
# rows 11-14 (Ringworm/Warts Molluscum/Serorrheic Keratoses/Nail Fungus stubs) up
# to rows 10-13, matching the final layout.
$ws.Rows.Item(10).Delete()

$v1 = @"
Eggplant and apple cider vinegar.
Baking soda and coconut oil paste.
Black salve or bloodroot.
Oils: black raspberry seed, frankincense, and myrrh.
Iodine.
Vitamin C.
"@
$ws.Range("D6").Value = $v1
$ws.Range("D6").WrapText = $true
$ws.Range("D6").VerticalAlignment = -4160

$v2 = @"
Eggplant and apple cider vinegar.
Baking soda and coconut oil paste.
Black salve or bloodroot.
Oils: black raspberry seed, frankincense, and myrrh.
Iodine.
Vitamin C.
"@
$ws.Range("D7").Value = $v2
$ws.Range("D7").WrapText = $true
$ws.Range("D7").VerticalAlignment = -4160

$v3 = @"
Ringworm
"@
$ws.Range("A10").Value = $v3
$ws.Range("A10").WrapText = $true
$ws.Range("A10").VerticalAlignment = -4160

$v4 = @"
A mild solution of bleach and water may ease inflammation and itching, as well as killing the bacteria that can cause skin infections when you have eczema. ...
Add apple cider vinegar to bath water. ...
Use mild soaps or body cleansers. ...
Moisturize your skin twice a day. ...
Apply coconut oil to damp skin once or twice a day.
"@
$ws.Range("D9").Value = $v4
$ws.Range("D9").WrapText = $true
$ws.Range("D9").VerticalAlignment = -4160

$v5 = @"
Ringworm often causes a ring-shaped rash that is itchy, scaly and slightly raised. The rings usually start small and then expand outward. Ringworm of the body (tinea corporis) is a rash caused by a fungal infection. It's usually an itchy, circular rash with clearer skin in the middle
"@
$ws.Range("B10").Value = $v5
$ws.Range("B10").WrapText = $true
$ws.Range("B10").VerticalAlignment = -4160

$v6 = @"
Itchy skin.
Ring-shaped rash.
Red, scaly, cracked skin.
Hair loss
"@
$ws.Range("C10").Value = $v6
$ws.Range("C10").WrapText = $true
$ws.Range("C10").VerticalAlignment = -4160

$v7 = @"
Garlic. - Garlic paste may be used as a topical treatment, although no studies have been conducted on its use. ...
Soapy water. ...
Apple cider vinegar. ...
Aloe vera. ...
Coconut oil. ...
Grapefruit seed extract. ...
Turmeric. ...
Powdered licorice
"@
$ws.Range("D10").Value = $v7
$ws.Range("D10").WrapText = $true
$ws.Range("D10").VerticalAlignment = -4160

$v8 = @"
Warts Molluscum
"@
$ws.Range("A11").Value = $v8
$ws.Range("A11").WrapText = $true
$ws.Range("A11").VerticalAlignment = -4160

$v9 = @"
They often have a pearly appearance. They're usually smooth and firm. In most people, the lesions range from about the size of a pinhead to as large as a pencil eraser (2 to 5 millimeters in diameter). They may become itchy, sore, red, and/or swollen.
"@
$ws.Range("B11").Value = $v9
$ws.Range("B11").WrapText = $true
$ws.Range("B11").VerticalAlignment = -4160

$v10 = @"
First sign
The bumps appear on the skin between 2 and 8 weeks after you get the virus that causes this skin infection.
When the bumps first appear, you usually see ones that are small, firm, pink, flesh-colored, or white. These bumps will likely get bigger.
"@
$ws.Range("C11").Value = $v10
$ws.Range("C11").WrapText = $true
$ws.Range("C11").VerticalAlignment = -4160

$v11 = @"
Just apply one cup of cider vinegar to bath water and submerge nightly for 10 minutes or so. Cold compresses can also be applied to itchy or irritated bumps. How can it be prevented? Within reason, try to prevent skin-to-skin contact with someone who has molluscum
"@
$ws.Range("D11").Value = $v11
$ws.Range("D11").WrapText = $true
$ws.Range("D11").VerticalAlignment = -4160

$v12 = @"
Serorrheic Keratoses
"@
$ws.Range("A12").Value = $v12
$ws.Range("A12").WrapText = $true
$ws.Range("A12").VerticalAlignment = -4160

$v13 = @"
A seborrheic keratosis (seb-o-REE-ik ker-uh-TOE-sis) is a common noncancerous (benign) skin growth. People tend to get more of them as they get older. Seborrheic keratoses are usually brown, black or light tan. The growths (lesions) look waxy or scaly and slightly raised.
"@
$ws.Range("B12").Value = $v13
$ws.Range("B12").WrapText = $true
$ws.Range("B12").VerticalAlignment = -4160

$v14 = @"
Be slightly raised from the surrounding skin.
Be white or light tan in appearance, which may darken to brown or black.
Have a waxy, pasted-on look.
Look scaly or like a wart.
Be well-defined from the surrounding skin.
Not usually cause pain but may sometimes itch a little.
"@
$ws.Range("C12").Value = $v14
$ws.Range("C12").WrapText = $true
$ws.Range("C12").VerticalAlignment = -4160

$v15 = @"
There are no proven home remedies for seborrheic keratoses. Lemon juice or vinegar can irritate the skin, possibly causing the lesion to dry and crumble. However, there is no evidence that this is safe or effective.
"@
$ws.Range("D12").Value = $v15
$ws.Range("D12").WrapText = $true
$ws.Range("D12").VerticalAlignment = -4160

$v16 = @"
Nail fungus is a common infection of the nail. It begins as a white or yellow-brown spot under the tip of your fingernail or toenail. As the fungal infection goes deeper, the nail may discolor, thicken and crumble at the edge. Nail fungus can affect several nails
"@
$ws.Range("B13").Value = $v16
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = -4160

$v17 = @"
Diabetes and think you're developing nail fungus.
Bleeding around the nails.
Swelling or pain around the nails.
Difficulty walking
"@
$ws.Range("C13").Value = $v17
$ws.Range("C13").WrapText = $true
$ws.Range("C13").VerticalAlignment = -4160

$v18 = @"
Oils. Some oils contain antifungal, anti-bacterial, and/or antiseptic qualities. ...
Garlic. Garlic includes strong antifungal properties, making this a top home remedy for treating particular strains of toenail fungus. ...
Snakeroot Extract. ...
Baking Soda. ...
Mentholated Topical Ointment.
"@
$ws.Range("D13").Value = $v18
$ws.Range("D13").WrapText = $true
$ws.Range("D13").VerticalAlignment = -4160

$v19 = @"
No home treatment, Immediately meet a doctor
"@
$ws.Range("D8").Value = $v19
$ws.Range("D8").WrapText = $true
$ws.Range("D8").VerticalAlignment = -4160

# Row heights to fit the newly-added / expanded content
$ws.Rows.Item(5).RowHeight = 160
$ws.Rows.Item(6).RowHeight = 110
$ws.Rows.Item(7).RowHeight = 166
$ws.Rows.Item(10).RowHeight = 111
$ws.Rows.Item(11).RowHeight = 108
$ws.Rows.Item(12).RowHeight = 97
$ws.Rows.Item(13).RowHeight = 119

# Restore the selection Excel shows after this edit
$ws.Range("D7").Select()

